$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 142941.28
$ws.Range("I9").Value = 250072.25
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 250072.25
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -249903.25
$ws.Range("N9").Value = -438
$ws.Range("H17").Value = 3392.25
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3392.25
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 10176.75
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -10512.75
$ws.Range("H21").Value = 6102.2856
$ws.Range("I21").Value = 3786
$ws.Range("K21").Value = 3786
$ws.Range("M21").Value = -3318
$ws.Range("H23").Value = 6102.2856
$ws.Range("I23").Value = 3786
$ws.Range("K23").Value = 3786
$ws.Range("M23").Value = -3552
$ws.Range("H28").Value = 5250
$ws.Range("J28").Value = 500
$ws.Range("L28").Value = 500
$ws.Range("N28").Value = -1470
$ws.Range("H33").Value = 356.33334
$ws.Range("I33").Value = 356.33334
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 356.33334
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -127.33334
$ws.Range("N33").ClearContents()
$ws.Range("H64").Value = 5799.5
$ws.Range("I64").Value = 5800
$ws.Range("J64").Value = 5799
$ws.Range("K64").Value = 5800
$ws.Range("L64").Value = 5799
$ws.Range("M64").Value = -5552
$ws.Range("N64").Value = -6295
$ws.Range("H67").Value = 5799.5
$ws.Range("I67").Value = 5800
$ws.Range("J67").Value = 5799
$ws.Range("K67").Value = 5800
$ws.Range("L67").Value = 5799
$ws.Range("M67").Value = -4942
$ws.Range("N67").Value = -7515
$ws.Range("H74").Value = 4190.7334
$ws.Range("I74").Value = 4190.7334
$ws.Range("K74").Value = 4190.7334
$ws.Range("M74").Value = -3254.7334
$ws.Range("H77").Value = 4190.7334
$ws.Range("I77").Value = 4190.7334
$ws.Range("K77").Value = 20953.667
$ws.Range("M77").Value = -16273.667
$ws.Range("H129").Value = 1426.2727
$ws.Range("I129").Value = 1299
$ws.Range("J129").Value = 1499
$ws.Range("K129").Value = 3897
$ws.Range("L129").Value = 4497
$ws.Range("M129").Value = 1103
$ws.Range("N129").Value = -14497

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5409.55
$ws.Range("I32").Value = 5355.974
$ws.Range("K32").Value = 5355.974
$ws.Range("M32").Value = -5068.974
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H74").Value = 4976.6177
$ws.Range("I74").Value = 4669.4136
$ws.Range("K74").Value = 4669.4136
$ws.Range("M74").Value = -3795.4136
$ws.Range("H77").Value = 4976.6177
$ws.Range("I77").Value = 4669.4136
$ws.Range("K77").Value = 23347.068
$ws.Range("M77").Value = -18979.068
$ws.Range("H97").Value = 925.4286
$ws.Range("I97").Value = 413.16666
$ws.Range("K97").Value = 413.16666
$ws.Range("M97").Value = 82.83334000000002
$ws.Range("H132").Value = 3223.3333
$ws.Range("I132").Value = 2335.5
$ws.Range("K132").Value = 7006.5
$ws.Range("M132").Value = -4476.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 387.25
$ws.Range("I80").Value = 117
$ws.Range("J80").Value = 522.375
$ws.Range("K80").Value = 117
$ws.Range("L80").Value = 522.375
$ws.Range("M80").Value = 881
$ws.Range("N80").Value = -2518.375
$ws.Range("H83").Value = 387.25
$ws.Range("I83").Value = 117
$ws.Range("J83").Value = 522.375
$ws.Range("K83").Value = 585
$ws.Range("L83").Value = 2611.875
$ws.Range("M83").Value = 4407
$ws.Range("N83").Value = -12595.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 35000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H14").Value = 500.22223
$ws.Range("I14").Value = 500.22223
$ws.Range("K14").Value = 1500.66669
$ws.Range("M14").Value = -1327.66669
$ws.Range("H33").Value = 80
$ws.Range("I33").Value = 80
$ws.Range("K33").Value = 480
$ws.Range("M33").Value = -197

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2266.6667
$ws.Range("I102").Value = 2266.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2266.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -644.6667000000002
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 749.1667
$ws.Range("I7").Value = 499.33334
$ws.Range("J7").Value = 999
$ws.Range("K7").Value = 499.33334
$ws.Range("L7").Value = 999
$ws.Range("M7").Value = -387.33334
$ws.Range("N7").Value = -1223
$ws.Range("H22").Value = 11057.917
$ws.Range("I22").Value = 23232
$ws.Range("K22").Value = 23232
$ws.Range("M22").Value = -22937
$ws.Range("H27").Value = 11057.917
$ws.Range("I27").Value = 23232
$ws.Range("K27").Value = 23232
$ws.Range("M27").Value = -23125
$ws.Range("H93").Value = 888.1177
$ws.Range("I93").Value = 881.7273
$ws.Range("K93").Value = 881.7273
$ws.Range("M93").Value = 366.2727
$ws.Range("H122").Value = 3442.4443
$ws.Range("I122").Value = 2997.1667
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 8991.500100000001
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -6541.500100000001
$ws.Range("N122").Value = -17899
$ws.Range("H126").Value = 749.1667
$ws.Range("I126").Value = 499.33334
$ws.Range("J126").Value = 999
$ws.Range("K126").Value = 1498.00002
$ws.Range("L126").Value = 2997
$ws.Range("M126").Value = 971.9999800000001
$ws.Range("N126").Value = -7937
$ws.Range("H132").Value = 4657.8
$ws.Range("I132").Value = 1606.75
$ws.Range("K132").Value = 4820.25
$ws.Range("M132").Value = -2290.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 22514
$ws.Range("I54").Value = 21000
$ws.Range("J54").Value = 22892.5
$ws.Range("K54").Value = 21000
$ws.Range("L54").Value = 22892.5
$ws.Range("M54").Value = -20480
$ws.Range("N54").Value = -23932.5
$ws.Range("H132").Value = 2185.2727
$ws.Range("I132").Value = 691.7143
$ws.Range("K132").Value = 2075.1429
$ws.Range("M132").Value = 454.8571000000002
